# v3.0 update FCI 27/1/2023
# Adds a new date column (C) with this week's NAV figures, and reorders the
# fund rows so the individual funds come first, followed by avg and total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header in column C, matching the style already used by B1 ---
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row labels (column A), reordered: funds, then avg, then total ---
$labels = @(
    "Alpha Acciones",
    "Delta Recursos Naturales",
    "Delta Select",
    "Delta gestion V",
    "HF Acciones Lideres",
    "Pionero Acciones",
    "avg",
    "total"
)

# --- Existing (06-01-2023) values, column B, aligned to the new row order ---
$colB = @(
    77443.96000000001,
    571463.61,
    203026.25,
    15099.4,
    69528.98,
    51312.12,
    164645.72,
    987874.3199999999
)

# --- New (13-01-2023) values, column C, aligned to the new row order ---
$colC = @(
    80682.92999999999,
    574928.9300000001,
    202954.79,
    24715.09,
    70357.17999999999,
    47973.37,
    166935.38,
    1001612.29
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $labels[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}
